# Request Science Gateway Community Accounts v1 - Google Drive -> GitHub/RTD conversion
#
# 1) Re-point the "Science Gateways Description" hyperlink at the same
#    target (it gets renumbered to rId6 once the unused comments
#    relationship is dropped from the package).
# 2) Resolve/remove the two open review comments ("Point to allocations
#    website ..." and "where is this process documentation ...") and
#    fold their commented spans back into plain text.
# 3) Reword the "request community accounts ..." sentence to describe
#    the individual-ticket-request process instead of the old
#    registration-process language that the comments were attached to.
# 4) Tidy up the leftover double space in "accessible.  It may take ..."
#    now that the comment reference between the two runs is gone.

$d = $word.ActiveDocument

# --- 1) hyperlink -----------------------------------------------------
$h = $d.Hyperlinks(1)
$h.Address = "https://docs.google.com/document/d/10shxQ6PallP6EIzFyL8_08cbKzzDvcqg_gVS6-Iv4lQ/"

# --- 2) drop the two John-Paul Navarro review comments -----------------
while ($d.Comments.Count -gt 0) {
    $d.Comments(1).Delete()
}

# --- 3) reword the community-accounts-request sentence ----------------
$rsq = [char]8217
$oldSentence = "Community accounts are UNIX accounts on ACCESS resource providers that can be used by a science gateway" + $rsq + "s users through the science gateway" + $rsq + "s user or programming interface. Science gateway providers request community accounts as part of the registration process. Accounts are created by ACCESS internal mechanisms. "
$newSentence = "Community accounts are UNIX accounts on ACCESS resource providers that can be used by a science gateway" + $rsq + "s users through the science gateway" + $rsq + "s user or programming interface. Science gateway providers request community accounts in individual ticket requests directed to the ACCESS Resources Providers, the ACCESS Integration Coordinator assigned can help make these requsts. "

$rng = $d.Content
$rng.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2)

# --- 4) single space before "It may take hours..." ---------------------
$rng2 = $d.Content
$rng2.Find.Execute("  It may take hours to days for accounts to be created after registration.", $true, $false, $false, $false, $false, $true, 1, $false, " It may take hours to days for accounts to be created after registration.", 2)
